# Update the three Job# values in MDSiTestResult.xlsx to the newest
# batch of job numbers (config-file refresh). The cells are stored as
# text (shared strings), so the cell is temporarily formatted as Text
# before the new value is typed in, then the explicit formatting is
# cleared again so the cell keeps using the sheet's default style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B4").NumberFormat = "@"

$ws.Range("B2").Value = "32255213"
$ws.Range("B3").Value = "32255214"
$ws.Range("B4").Value = "32255215"

$ws.Range("B2:B4").Style = "Normal"
